$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-02-17"

# Update the row label for February to reflect the new "through" date
$ws.Range("A3").Value = "February (through 02-17)"

# Update January 2022 value (column I, row 2)
$ws.Range("I2").Value = 160

# Update February row (row 3) values for years 2015-2021 and 2022
$ws.Range("B3").Value = 6
$ws.Range("C3").Value = 23
$ws.Range("D3").Value = 39
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 17
$ws.Range("G3").Value = 42
$ws.Range("H3").Value = 77
$ws.Range("I3").Value = 78

# Update Total row (row 4) values for years 2015-2021 (2022 total stays 238)
$ws.Range("B4").Value = 32
$ws.Range("C4").Value = 74
$ws.Range("D4").Value = 114
$ws.Range("E4").Value = 120
$ws.Range("F4").Value = 66
$ws.Range("G4").Value = 116
$ws.Range("H4").Value = 294
